$d = $word.ActiveDocument

$sentence = ". The graph plan for this is displayed below. The red lines indicate mutexes and the labels next to each arc indicates the type of mutex connecting the two components. The graph plan was continued until the last two state layers s2 and s3 were the same and the graph was no longer evolving. The solution to this problem is :"

# ---------------------------------------------------------------------------
# 1. Fix the typo "located int the" -> "located in the" (Question 4 paragraph)
# ---------------------------------------------------------------------------
$f = $d.Content.Find
$f.Execute("located int the", $true, $false, $false, $false, $false, $true, 1, $false, "located in the", 2)

# ---------------------------------------------------------------------------
# 2. Remove the stray "_GoBack" bookmark that currently sits in the empty
#    paragraph right before "Question 4".
# ---------------------------------------------------------------------------
$oldBm = $d.Bookmarks("_GoBack")
$oldBm.Delete()

# NOTE: both paragraphs edited below share the identical ". The graph
# plan..." sentence, so edits are applied in reverse document order (the
# later Q2_Problem2.txt paragraph first) so earlier character offsets found
# via Find stay valid for the edit that follows.

# ---------------------------------------------------------------------------
# 3. Question 2 / Q2_Problem2.txt paragraph (appears second in the document):
#    "...described in the file Q2_Problem2.txt. The graph plan..."
#    becomes
#    "...described in the file Q2_Problem2.txt located in the Assign4
#     folder. The graph plan..."
# ---------------------------------------------------------------------------
$p2 = $d.Content
$p2.Find.Execute("Q2_Problem2.txt")
$afterP2 = $p2.End

$r3 = $d.Range($afterP2, $afterP2 + 400)
$r3.Find.Execute($sentence)
$newText2 = " located in the Assign4 folder"
$r3.Text = $newText2 + $r3.Text

$r3b = $d.Range($afterP2, $afterP2 + $newText2.Length)
$r3b.Font.Size = 12
$r3b.Font.Size = 11

# ---------------------------------------------------------------------------
# 4. Question 2 / Q2_Problem1.txt paragraph (appears first in the document):
#    "...described in the file Q2_Problem1.txt. The graph plan..."
#    becomes
#    "...described in the file Q2_Problem1.txt located in the Assign4 folder.
#     The graph plan..." with a new _GoBack bookmark right after the
#     inserted text.
# ---------------------------------------------------------------------------
$p1 = $d.Content
$p1.Find.Execute("Q2_Problem1.txt")
$p1.Text = $p1.Text + " "
$insertPos = $p1.End

$newText1 = "located in the Assign4 folder"
$r1 = $d.Range($insertPos, $insertPos + 400)
$r1.Find.Execute($sentence)
$r1.Text = $newText1 + $r1.Text

# Force the newly inserted text to live in its own run (it currently shares
# the run with the ". The graph plan..." text because both ended up with the
# exact same run properties) by nudging a formatting property off and back.
$r1b = $d.Range($insertPos, $insertPos + $newText1.Length)
$r1b.Font.Size = 12
$r1b.Font.Size = 11

# Re-insert the "_GoBack" bookmark right after the newly inserted text.
$bmRange = $d.Range($insertPos + $newText1.Length, $insertPos + $newText1.Length)
$d.Bookmarks.Add("_GoBack", $bmRange)
